$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: Jasmine
$ws.Range("A14").Value = "Jasmine"
$ws.Range("B14").Value = "2.3.4"
$ws.Range("C14").Value = "https://github.com/jasmine/jasmine/blob/master/MIT.LICENSE"
$ws.Range("D14").Value = "http://jasmine.github.io/2.3/introduction.html"
$ws.Range("E14").Value = "MIT Licensed"
$ws.Range("F14").Value = "JavaScript Testing"

# New row 15: Karma
$ws.Range("A15").Value = "Karma"
$ws.Range("B15").Value = "0.12.36"
$ws.Range("C15").Value = "https://github.com/karma-runner/karma/blob/master/LICENSE"
$ws.Range("D15").Value = "http://karma-runner.github.io/0.12/index.html"
$ws.Range("E15").Value = "MIT Licensed"
$ws.Range("F15").Value = "JavaScript Testing"

# New row 16: Karma-Jasmine
$ws.Range("A16").Value = "Karma-Jasmine"
$ws.Range("B16").Value = "0.3.5"
$ws.Range("C16").Value = "https://github.com/karma-runner/karma-jasmine/blob/master/LICENSE"
$ws.Range("D16").Value = "https://github.com/karma-runner/karma-jasmine"
$ws.Range("E16").Value = "MIT Licensed"
$ws.Range("F16").Value = "JavaScript Testing"

# New row 17: Karma-Junit-Reporter
$ws.Range("A17").Value = "Karma-Junit-Reporter"
$ws.Range("B17").Value = "0.2.2"
$ws.Range("C17").Value = "https://github.com/karma-runner/karma-junit-reporter/blob/master/LICENSE"
$ws.Range("D17").Value = "https://github.com/karma-runner/karma-junit-reporter"
$ws.Range("E17").Value = "MIT Licensed"
$ws.Range("F17").Value = "JavaScript Testing"

# New row 18: Karma-Chrome-Launcher
$ws.Range("A18").Value = "Karma-Chrome-Launcher"
$ws.Range("B18").Value = "0.1"
$ws.Range("C18").Value = "https://github.com/karma-runner/karma-chrome-launcher/blob/master/LICENSE"
$ws.Range("D18").Value = "https://github.com/karma-runner/karma-chrome-launcher"
$ws.Range("E18").Value = "MIT Licensed"
$ws.Range("F18").Value = "Starts Google Chrome"

# Resize the table to include the newly added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K29"))

# Move the selection to match the post-edit state
$ws.Range("F17").Select()
